$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58, shifting existing rows 58-60 down to 59-61
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the new weekly data point
$ws.Cells.Item(58, 1).Value = 7
$ws.Cells.Item(58, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(58, 3).Value = "Ñuble"
$ws.Cells.Item(58, 4).Value = 44595
$ws.Cells.Item(58, 5).Value = 16
$ws.Cells.Item(58, 6).Value = 100112022
$ws.Cells.Item(58, 7).Value = "Arveja Verde"
$ws.Cells.Item(58, 8).Value = "Sin especificar"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 60
$ws.Cells.Item(58, 11).Value = 23000
$ws.Cells.Item(58, 12).Value = 24000
$ws.Cells.Item(58, 13).Value = 23500
$ws.Cells.Item(58, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(58, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(58, 16).Value = 940
$ws.Cells.Item(58, 17).Value = 25
$ws.Cells.Item(58, 18).Value = "Hortaliza"
